$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cells (bold, centered, thin border) by copying H1's format.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I0/IF data for rows 2-74 (row, I0, IF)
$data = @(
    "2,9,9",
    "3,9,9",
    "4,9,9",
    "5,9,9",
    "6,9,9",
    "7,9,9",
    "8,9,9",
    "9,9,9",
    "10,9,9",
    "11,8,9",
    "12,9,9",
    "13,9,9",
    "14,9,9",
    "15,8,9",
    "16,8,9",
    "17,9,9",
    "18,9,9",
    "19,9,10",
    "20,9,9",
    "21,9,9",
    "22,9,9",
    "23,9,9",
    "24,9,9",
    "25,9,9",
    "26,9,9",
    "27,9,9",
    "28,9,9",
    "29,8,8",
    "30,9,9",
    "31,9,10",
    "32,9,9",
    "33,9,9",
    "34,8,8",
    "35,9,9",
    "36,8,9",
    "37,7,7",
    "38,9,9",
    "39,9,9",
    "40,9,9",
    "41,9,9",
    "42,9,9",
    "43,9,9",
    "44,8,8",
    "45,9,9",
    "46,10,10",
    "47,9,9",
    "48,6,6",
    "49,6,6",
    "50,10,10",
    "51,10,10",
    "52,9,9",
    "53,8,8",
    "54,10,10",
    "55,9,9",
    "56,9,9",
    "57,8,8",
    "58,9,9",
    "59,8,9",
    "60,9,9",
    "61,9,9",
    "62,9,9",
    "63,9,9",
    "64,9,9",
    "65,9,9",
    "66,9,9",
    "67,8,8",
    "68,7,8",
    "69,9,9",
    "70,6,6",
    "71,4,4",
    "72,6,6",
    "73,6,6",
    "74,5,5"
)

foreach ($entry in $data) {
    $parts = $entry.Split(",")
    $r = [int]$parts[0]
    $i0 = [int]$parts[1]
    $iF = [int]$parts[2]
    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $iF
}
